# Generate Report for Handback
# Refresh the handback-status timestamps to reflect the latest report run.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-21 05:10:01"

# zh-cn sheet: handoff/handback datetimes for the first file
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-21 05:09:56"
$wsZhCn.Range("K2").Value = "2016-08-21 05:10:27"

# de-de sheet: handback datetime for the first file
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-21 05:10:35"
